$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.000009318123435519965
$ws.Range("C2").Value = 0.05231270169004087
$ws.Range("D2").Value = 2938.103010863317
$ws.Range("E2").Value = 198602002.3250627
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 198604940.4803956

$ws.Range("B3").Value = 3.182878228561681
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 3.082599426703578
$ws.Range("E3").Value = 6.48142807727062
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 14.40014219143469

$ws.Range("B4").Value = 1.505614041169197
$ws.Range("C4").Value = 0.05231270169004087
$ws.Range("D4").Value = 0.7127328510149897
$ws.Range("E4").Value = 0.4998867070740569
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 2.770546300948285

$ws.Range("B5").Value = 0.7287194209349384
$ws.Range("C5").Value = 9.226618575922256
$ws.Range("D5").Value = 3.082599426703578
$ws.Range("E5").Value = 6.48142807727062
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 19.51936550083139

$ws.Range("B6").Value = 1.505614041169197
$ws.Range("C6").Value = 1.65323645889881
$ws.Range("D6").Value = 0.7127328510149897
$ws.Range("E6").Value = 246.9852506941017
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 250.8568340451847
